# Fix bug in grade structure:
#  - Rename headers GK1/CK1/Bonu -> GK/CK/Bonus
#  - Re-point row 2 to "Hoa Pham 1" with new grades, row 3 to "Hoa Pham 2"
#  - Add a new row 4 for "Hoa Pham PVH 1"
#
# Numeric-looking grade values must stay stored as TEXT (matching the
# workbook's existing "numbers stored as text" convention), so each such
# cell is written while the cell is text-formatted, then the formatting
# is cleared again so no stray number-format/style change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Header row
$ws.Range("B1").Value = "GK"
$ws.Range("C1").Value = "CK"
$ws.Range("D1").Value = "Bonus"

# Row 2 -> "Hoa Pham 1"
$ws.Range("A2").Value = "Hoa Pham 1"
Set-TextValue "B2" "20"
Set-TextValue "C2" "50"
Set-TextValue "D2" "10"
Set-TextValue "E2" "80"

# Row 3 -> "Hoa Pham 2"
$ws.Range("A3").Value = "Hoa Pham 2"
Set-TextValue "B3" "10"
Set-TextValue "C3" "50"
Set-TextValue "D3" "10"
Set-TextValue "E3" "70"

# New row 4 -> "Hoa Pham PVH 1"
$ws.Range("A4").Value = "Hoa Pham PVH 1"
Set-TextValue "B4" "5"
Set-TextValue "C4" "50"
Set-TextValue "D4" "10"
Set-TextValue "E4" "65"
